$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Status of the first item moved from "En proceso" to "Cerrada"
$ws.Range("F4").Value = "Cerrada"

# Fill in the actual closing date, matching the planned closing date
$ws.Range("E4").Value = $ws.Range("D4").Value2

# Update the selected/active cell to F5
$ws.Range("F5").Select()
